$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the username (email) and password values stored in A2 / B2.
$ws.Range("A2").Value = "ravi.ranjan@onpassive.com"
$ws.Range("B2").Value = "Onpassive90@"
